$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix training data issue: the "Date" column (BF) held the source file
# name ("6-5-2007-08") instead of the actual game date. NBA stats for a
# game played on 2008-06-05 were off by a day due to how the stats were
# shown, so correct the Date column to the ISO date string "2008-06-05"
# for every data row (rows 2-31).
#
# Force the range to text formatting first so Excel stores the literal
# string "2008-06-05" instead of auto-converting it to a date serial
# number, then restore the cell style so no extra formatting is applied.
$dateRange = $ws.Range("BF2:BF31")
$dateRange.NumberFormat = "@"
for ($row = 2; $row -le 31; $row++) {
    $ws.Range("BF$row").Value = "2008-06-05"
}
$dateRange.Style = "Normal"
